# Rearrange the parks alphabetically.
#
# Each worksheet lists model-summary rows for predictors including the
# "Park*" dummy variables. Rows 5-11 hold (in this order, before the edit):
#   5  ParkHainich
#   6  ParkHunsrueck
#   7  ParkSaechs_Schw
#   8  ParkJasmund
#   9  ParkKellerwald
#   10 ParkEifel
#   11 ParkVorpomm
#
# The park labels are re-sorted alphabetically (Eifel, Hainich, Hunsrueck,
# Jasmund, Kellerwald, Saechs_Schw, Vorpomm), which only moves rows 5, 6, 7
# and 10 (Jasmund/Kellerwald/Vorpomm were already in the right spot). The
# label cell (column A) of every row 5-11 keeps its row, but the model
# statistics that travel with each label move together with it, and on the
# "API" sheet (the only one whose model failed to converge for some of
# these levels) the refit drops the numbers for three of the levels.

$wb = $excel.ActiveWorkbook

function Set-ParkRow {
    param(
        $ws,
        [int]$row,
        [string]$label,
        $bVal,
        $cVal,
        $dVal
    )
    $ws.Cells.Item($row, 1).Value = $label
    if ($null -eq $bVal) {
        $ws.Cells.Item($row, 2).Value = $null
    } else {
        $ws.Cells.Item($row, 2).Value = $bVal
    }
    $ws.Cells.Item($row, 3).Value = $cVal
    if ($null -eq $dVal) {
        $ws.Cells.Item($row, 4).Value = $null
    } else {
        $ws.Cells.Item($row, 4).Value = $dVal
    }
}

# --- Sheet 1: API ---------------------------------------------------------
# This sheet's glm had separation issues; the refit (with the new park
# reference level) changes B2/C2 too, and drops B/C/D for Hunsrueck,
# Jasmund, Saechs_Schw and Vorpomm (now rows 7, 8, 10 and 11).
$ws = $wb.Worksheets.Item("API")
$ws.Cells.Item(2, 2).Value = 0.2467
$ws.Cells.Item(2, 3).Value = 85816.1197

Set-ParkRow $ws 5  "ParkEifel"       1.7447    85816.1197 1
Set-ParkRow $ws 6  "ParkHainich"     1.8242    85816.1197 1
Set-ParkRow $ws 7  "ParkHunsrueck"   $null     0          $null
Set-ParkRow $ws 8  "ParkJasmund"     $null     0          $null
Set-ParkRow $ws 9  "ParkKellerwald"  1.3494    85816.1197 1
Set-ParkRow $ws 10 "ParkSaechs_Schw" $null     0          $null
Set-ParkRow $ws 11 "ParkVorpomm"     $null     0          $null

# --- Remaining sheets: simple numeric model tables ------------------------
$sheets = @{
    "Anthropogenic pollution" = @{
        5  = @("ParkEifel",      -0.2778, 0.1692, 0.1006)
        6  = @("ParkHainich",     0.2237, 0.1549, 0.1486)
        7  = @("ParkHunsrueck",  -0.5301, 0.2023, 0.0088)
        10 = @("ParkSaechs_Schw", 0.241,  0.147,  0.1011)
    }
    "Industrial chemical" = @{
        5  = @("ParkEifel",       0.7964, 0.4608, 0.0839)
        6  = @("ParkHainich",    -4.1308, 31132.3064, 0.9999)
        7  = @("ParkHunsrueck",  -0.1033, 0.5976, 0.8627)
        10 = @("ParkSaechs_Schw", 0.3199, 0.4924, 0.516)
    }
    "PAH" = @{
        5  = @("ParkEifel",       0.0062, 0.0808, 0.9385)
        6  = @("ParkHainich",    -0.0508, 0.09,   0.5723)
        7  = @("ParkHunsrueck",   0.0425, 0.0822, 0.6055)
        10 = @("ParkSaechs_Schw",-0.1305, 0.0857, 0.128)
    }
    "PCP" = @{
        5  = @("ParkEifel",       0.5066, 0.4029, 0.2086)
        6  = @("ParkHainich",     0.2732, 0.4071, 0.5022)
        7  = @("ParkHunsrueck",   0.3417, 0.4097, 0.4043)
        10 = @("ParkSaechs_Schw", 0.0394, 0.4375, 0.9282)
    }
    "POP" = @{
        5  = @("ParkEifel",       0.2152, 0.1585, 0.1745)
        6  = @("ParkHainich",    -0.3065, 0.2529, 0.2254)
        7  = @("ParkHunsrueck",   0.2063, 0.1625, 0.2041)
        10 = @("ParkSaechs_Schw", 0.061,  0.1622, 0.707)
    }
    "Pesticide" = @{
        5  = @("ParkEifel",       0.1365, 0.6849, 0.842)
        6  = @("ParkHainich",     0.607,  0.5845, 0.299)
        7  = @("ParkHunsrueck",  -0.2895, 0.7234, 0.6891)
        10 = @("ParkSaechs_Schw",-2.0377, 2189.6613, 0.9993)
    }
    "Plasticizer" = @{
        5  = @("ParkEifel",      -0.0263, 0.2367, 0.9114)
        6  = @("ParkHainich",    -0.1329, 0.2392, 0.5783)
        7  = @("ParkHunsrueck",   0.0402, 0.2417, 0.8679)
        10 = @("ParkSaechs_Schw",-0.7947, 0.2635, 0.0026)
    }
}

foreach ($sheetName in $sheets.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $sheets[$sheetName]
    foreach ($row in $rows.Keys) {
        $vals = $rows[$row]
        Set-ParkRow $ws $row $vals[0] $vals[1] $vals[2] $vals[3]
    }
}
